# Alteração nos rótulos da tabela para transformar a primeira linha em
# cabeçalho automaticamente no Power BI.
# Prefixa os rótulos de ano / intervalo da linha 1 de cada planilha.

$wb = $excel.ActiveWorkbook

# Planilhas cujo cabeçalho é "Ano <ano>" nas colunas B:E
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Planilha cujo cabeçalho é "Intervalo <intervalo>" nas colunas B:E
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws.Range("B1").Value = "Intervalo 2015"
$ws.Range("C1").Value = "Intervalo 2015-2030"
$ws.Range("D1").Value = "Intervalo 2031-2040"
$ws.Range("E1").Value = "Intervalo 2041-2050"

# Planilha com apenas a coluna B no cabeçalho
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws.Range("B1").Value = "Ano 2015"
